$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-08-25T17:07:23"
$ws.Range("U4").Value = 44.87
$ws.Range("V4").Value = 40
$ws.Range("W4").Value = 36.89
$ws.Range("X4").Value = 34.38
$ws.Range("Y4").Value = 31.76
$ws.Range("Z4").Value = 34.48
$ws.Range("U6").Value = -2.11
$ws.Range("V6").Value = -1.8
$ws.Range("W6").Value = -1.96
$ws.Range("X6").Value = -1.51
$ws.Range("Y6").Value = -1.05
$ws.Range("Z6").Value = -1.07
$ws.Range("U8").Value = 0
$ws.Range("V8").Value = 0
$ws.Range("U9").Value = 44.82
$ws.Range("V9").Value = 44.31
$ws.Range("W9").Value = 36.85
$ws.Range("X9").Value = 35.22
$ws.Range("Y9").Value = 32.81
$ws.Range("Z9").Value = 35.59
$ws.Range("U11").Value = -2.15
$ws.Range("V11").Value = -0.78
$ws.Range("W11").Value = -1.99
$ws.Range("X11").Value = -0.67
$ws.Range("Y11").Value = 0
$ws.Range("Z11").Value = 0.04
$ws.Range("V12").Value = 3.29
$ws.Range("U13").Value = 0
$ws.Range("V13").Value = 0
$ws.Range("U14").Value = 44.82
$ws.Range("V14").Value = 44.31
$ws.Range("W14").Value = 36.85
$ws.Range("X14").Value = 35.22
$ws.Range("Y14").Value = 32.81
$ws.Range("Z14").Value = 35.62
$ws.Range("U16").Value = -2.15
$ws.Range("V16").Value = -0.78
$ws.Range("W16").Value = -1.99
$ws.Range("X16").Value = -0.67
$ws.Range("Y16").Value = 0
$ws.Range("Z16").Value = 0.07000000000000001
$ws.Range("V17").Value = 3.29
$ws.Range("U18").Value = 0
$ws.Range("V18").Value = 0
$ws.Range("V19").Value = 39.81
$ws.Range("W19").Value = 36.68
$ws.Range("X19").Value = 34.28
$ws.Range("Y19").Value = 31.79
$ws.Range("Z19").Value = 34.42
$ws.Range("U21").Value = -2.36
$ws.Range("V21").Value = -1.99
$ws.Range("W21").Value = -2.16
$ws.Range("X21").Value = -1.61
$ws.Range("Y21").Value = -1.02
$ws.Range("Z21").Value = -1.14
$ws.Range("U23").Value = 0
$ws.Range("V23").Value = 0
$ws.Range("V24").Value = 39.81
$ws.Range("W24").Value = 36.68
$ws.Range("X24").Value = 34.28
$ws.Range("Y24").Value = 31.79
$ws.Range("Z24").Value = 34.42
$ws.Range("U26").Value = -2.36
$ws.Range("V26").Value = -1.99
$ws.Range("W26").Value = -2.16
$ws.Range("X26").Value = -1.61
$ws.Range("Y26").Value = -1.02
$ws.Range("Z26").Value = -1.14
$ws.Range("U28").Value = 0
$ws.Range("V28").Value = 0
$ws.Range("U29").Value = 44.31
$ws.Range("V29").Value = 39.55
$ws.Range("W29").Value = 36.44
$ws.Range("X29").Value = 34.15
$ws.Range("Y29").Value = 31.76
$ws.Range("Z29").Value = 34.32
$ws.Range("U31").Value = -2.66
$ws.Range("V31").Value = -2.25
$ws.Range("W31").Value = -2.4
$ws.Range("X31").Value = -1.74
$ws.Range("Y31").Value = -1.05
$ws.Range("Z31").Value = -1.24
$ws.Range("U33").Value = 0
$ws.Range("V33").Value = 0
$ws.Range("U34").Value = 45.34
$ws.Range("V34").Value = 42.09
$ws.Range("W34").Value = 37.14
$ws.Range("X34").Value = 35.93
$ws.Range("Y34").Value = 33.55
$ws.Range("Z34").Value = 36.58
$ws.Range("U36").Value = -1.63
$ws.Range("V36").Value = 0.29
$ws.Range("W36").Value = -1.71
$ws.Range("X36").Value = 0.04
$ws.Range("Y36").Value = 0.74
$ws.Range("Z36").Value = 1.02
$ws.Range("U38").Value = 0
$ws.Range("V38").Value = 0
$ws.Range("U39").Value = 44.87
$ws.Range("V39").Value = 40
$ws.Range("W39").Value = 36.89
$ws.Range("X39").Value = 34.38
$ws.Range("Y39").Value = 31.76
$ws.Range("Z39").Value = 34.48
$ws.Range("U41").Value = -2.11
$ws.Range("V41").Value = -1.8
$ws.Range("W41").Value = -1.96
$ws.Range("X41").Value = -1.51
$ws.Range("Y41").Value = -1.05
$ws.Range("Z41").Value = -1.07
$ws.Range("U43").Value = 0
$ws.Range("V43").Value = 0
$ws.Range("U44").Value = 46.55
$ws.Range("V44").Value = 41.51
$ws.Range("W44").Value = 38.42
$ws.Range("X44").Value = 35.5
$ws.Range("Y44").Value = 32.58
$ws.Range("Z44").Value = 35.62
$ws.Range("U46").Value = -0.42
$ws.Range("V46").Value = -0.29
$ws.Range("W46").Value = -0.42
$ws.Range("X46").Value = -0.39
$ws.Range("Y46").Value = -0.23
$ws.Range("Z46").Value = 0.07000000000000001
$ws.Range("U48").Value = 0
$ws.Range("V48").Value = 0
$ws.Range("U49").Value = 48.63
$ws.Range("V49").Value = 43.95
$ws.Range("W49").Value = 40.63
$ws.Range("X49").Value = 36.4
$ws.Range("Z49").Value = 36.24
$ws.Range("U51").Value = 1.65
$ws.Range("V51").Value = 2.15
$ws.Range("W51").Value = 1.79
$ws.Range("X51").Value = 0.51
$ws.Range("U53").Value = 0
$ws.Range("V53").Value = 0
$ws.Range("U54").Value = 46.46
$ws.Range("V54").Value = 41.72
$ws.Range("W54").Value = 39.28
$ws.Range("X54").Value = 36.89
$ws.Range("Y54").Value = 33.89
$ws.Range("Z54").Value = 36.96
$ws.Range("U56").Value = -0.51
$ws.Range("V56").Value = -0.08
$ws.Range("W56").Value = 0.43
$ws.Range("X56").Value = 1
$ws.Range("Y56").Value = 1.08
$ws.Range("Z56").Value = 1.4
$ws.Range("U58").Value = 0
$ws.Range("V58").Value = 0
$ws.Range("U59").Value = 48.38
$ws.Range("V59").Value = 43
$ws.Range("W59").Value = 39.88
$ws.Range("X59").Value = 36.77
$ws.Range("Y59").Value = 33.68
$ws.Range("Z59").Value = 36.84
$ws.Range("U61").Value = 1.4
$ws.Range("V61").Value = 1.2
$ws.Range("W61").Value = 1.04
$ws.Range("X61").Value = 0.88
$ws.Range("Y61").Value = 0.88
$ws.Range("Z61").Value = 1.29
$ws.Range("U63").Value = 0
$ws.Range("V63").Value = 0
$ws.Range("U64").Value = 49.24
$ws.Range("V64").Value = 43.77
$ws.Range("W64").Value = 40.59
$ws.Range("X64").Value = 37.35
$ws.Range("Y64").Value = 34.14
$ws.Range("Z64").Value = 37.42
$ws.Range("U66").Value = 2.26
$ws.Range("V66").Value = 1.97
$ws.Range("W66").Value = 1.75
$ws.Range("X66").Value = 1.46
$ws.Range("Y66").Value = 1.33
$ws.Range("Z66").Value = 1.87
$ws.Range("U68").Value = 0
$ws.Range("V68").Value = 0
$ws.Range("U69").Value = 49.76
$ws.Range("V69").Value = 44.23
$ws.Range("W69").Value = 41.06
$ws.Range("X69").Value = 37.82
$ws.Range("Y69").Value = 34.61
$ws.Range("Z69").Value = 37.9
$ws.Range("U71").Value = 2.79
$ws.Range("V71").Value = 2.43
$ws.Range("W71").Value = 2.22
$ws.Range("X71").Value = 1.93
$ws.Range("Y71").Value = 1.8
$ws.Range("Z71").Value = 2.35
$ws.Range("U73").Value = 0
$ws.Range("V73").Value = 0
$ws.Range("U74").Value = 48.23
$ws.Range("V74").Value = 43.09
$ws.Range("W74").Value = 39.96
$ws.Range("X74").Value = 36.74
$ws.Range("Y74").Value = 33.58
$ws.Range("Z74").Value = 36.73
$ws.Range("U76").Value = 1.25
$ws.Range("V76").Value = 1.29
$ws.Range("W76").Value = 1.12
$ws.Range("X76").Value = 0.84
$ws.Range("Y76").Value = 0.77
$ws.Range("Z76").Value = 1.18
$ws.Range("U78").Value = 0
$ws.Range("V78").Value = 0
$ws.Range("U79").Value = 48.57
$ws.Range("V79").Value = 43.32
$ws.Range("W79").Value = 40.23
$ws.Range("X79").Value = 36.96
$ws.Range("Y79").Value = 33.78
$ws.Range("Z79").Value = 36.95
$ws.Range("U81").Value = 1.59
$ws.Range("V81").Value = 1.52
$ws.Range("W81").Value = 1.39
$ws.Range("X81").Value = 1.07
$ws.Range("Y81").Value = 0.97
$ws.Range("Z81").Value = 1.4
$ws.Range("U83").Value = 0
$ws.Range("V83").Value = 0
$ws.Range("U84").Value = 47.45
$ws.Range("V84").Value = 42.65
$ws.Range("W84").Value = 40.17
$ws.Range("X84").Value = 37.62
$ws.Range("Y84").Value = 34.54
$ws.Range("Z84").Value = 37.54
$ws.Range("U86").Value = 0.47
$ws.Range("V86").Value = 0.85
$ws.Range("W86").Value = 1.33
$ws.Range("X86").Value = 1.73
$ws.Range("Y86").Value = 1.73
$ws.Range("Z86").Value = 1.99
$ws.Range("U88").Value = 0
$ws.Range("V88").Value = 0
$ws.Range("U89").Value = 44.31
$ws.Range("V89").Value = 39.51
$ws.Range("W89").Value = 36.44
$ws.Range("X89").Value = 34.15
$ws.Range("Y89").Value = 31.76
$ws.Range("Z89").Value = 34.32
$ws.Range("U91").Value = -2.66
$ws.Range("V91").Value = -2.29
$ws.Range("W91").Value = -2.4
$ws.Range("X91").Value = -1.74
$ws.Range("Y91").Value = -1.05
$ws.Range("Z91").Value = -1.24
$ws.Range("U93").Value = 0
$ws.Range("V93").Value = 0
